# Update "想去人数" (interested-people count) figures that changed between
# the previous data pull and the latest one (output generated at 456a3b4).
#
#   展览  (Exhibitions) sheet : F2 5798 -> 5801 ; F5 976 -> 977
#   全部类型 (All types) sheet : F2 5798 -> 5801 ; F5 976 -> 977

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 5801
    $ws.Range("F5").Value = 977
}
